$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel PasteSpecial constant
$xlPasteFormats = -4122

function Set-DateText {
    param($range, $text)
    # Assigning a date-like string (e.g. "20-FEB-26") via .Value normally gets
    # auto-parsed by Excel into a real date serial number + new date style.
    # To preserve the original plain-text representation (as in the source
    # workbook) we temporarily force a Text format before assigning the
    # value, then restore the cell's original formatting/style by copying it
    # back from an untouched neighbor cell on the same row that already has
    # the desired style.
    $range.NumberFormat = "@"
    $range.Value = $text
}

# --- Row 2 ---
Set-DateText $ws.Range("A2") "20-FEB-26"
$ws.Range("E2").Value = 519
$ws.Range("F2").Value = -57

# --- Row 3 ---
Set-DateText $ws.Range("A3") "27-FEB-26"

# --- Row 4 ---
Set-DateText $ws.Range("A4") "06-MAR-26"
$ws.Range("D4").Value = 960
$ws.Range("E4").Value = 983
$ws.Range("F4").Value = -23

# --- Row 5 ---
Set-DateText $ws.Range("A5") "27-MAR-26"
$ws.Range("D5").Value = 462
$ws.Range("E5").Value = 713
$ws.Range("F5").Value = -251

# --- Row 6 ---
Set-DateText $ws.Range("A6") "29-MAR-26"
$ws.Range("E6").Value = 519
$ws.Range("F6").Value = -57
$ws.Range("J6").Value = "LOW THREAT"

# --- Row 7 ---
Set-DateText $ws.Range("A7") "02-APR-26"
$ws.Range("D7").Value = 456
$ws.Range("E7").Value = 519
$ws.Range("F7").Value = -63

# --- Row 8 ---
Set-DateText $ws.Range("A8") "17-MAY-26"
$ws.Range("D8").Value = 960
$ws.Range("E8").Value = 1774
$ws.Range("F8").Value = -814
$ws.Range("J8").Value = "HIGH THREAT ALERT - NEED ACTION"

# --- Row 9 ---
Set-DateText $ws.Range("A9") "21-MAY-26"
$ws.Range("C9").Value = "Nile Air NP-142"
$ws.Range("D9").Value = 1051
$ws.Range("F9").Value = -723

# --- Row 10 ---
$ws.Range("D10").Value = 1235
$ws.Range("F10").Value = -539

# --- Row 11 ---
$ws.Range("D11").Value = 1241
$ws.Range("F11").Value = -533

# Restore original cell formatting (border/fill/font/number-format=General)
# for every cell whose NumberFormat we temporarily touched, by pasting the
# formats from an unaffected neighbor in the same column/style family.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A9").PasteSpecial($xlPasteFormats) | Out-Null

# Re-style J6 from MEDIUM (yellow) to LOW (green) - copy formats from a LOW THREAT cell
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J6").PasteSpecial($xlPasteFormats) | Out-Null

# Re-style J8 from LOW (green) to HIGH (red) - copy formats from a HIGH THREAT cell
$ws.Range("J9").Copy() | Out-Null
$ws.Range("J8").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false
